$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row {
    param($r, $values)
    for ($i = 0; $i -lt $values.Count; $i++) {
        $ws.Cells.Item($r, $i + 1).Value = $values[$i]
    }
}

# Row 2: ECs -> FAPs
Set-Row 2 @("ECs","Ccl25","Ccr10","FAPs",3,1,4.528161333333333,13.584484,0.3225352762763812,0.3225352762763812,2,0.6666666666666666,1.036595333333333,3.109786,0.393072250513715,0.393072250513715,4.693870906713777,42.244838160424,0.12677966691602,0.12677966691602)

# Row 3: ECs -> sCs
Set-Row 3 @("ECs","Ccl25","Ccr10","sCs",3,1,4.528161333333333,13.584484,0.3225352762763812,0.3225352762763812,3,1,1.600567,4.801701,0.6069277494862849,0.6069277494862849,7.247625600809333,65.228630407284,0.1957556093603612,0.1957556093603612)

# Row 4: FAPs -> FAPs
Set-Row 4 @("FAPs","Ccl25","Ccr10","FAPs",3,1,6.270951666666666,18.812855,0.4466720550425397,0.4466720550425397,2,0.6666666666666666,1.036595333333333,3.109786,0.393072250513715,0.393072250513715,6.500439233225555,58.50395309903,0.1755743899171571,0.1755743899171571)

# Row 5: FAPs -> sCs
Set-Row 5 @("FAPs","Ccl25","Ccr10","sCs",3,1,6.270951666666666,18.812855,0.4466720550425397,0.4466720550425397,3,1,1.600567,4.801701,0.6069277494862849,0.6069277494862849,10.03707829626167,90.333704666355,0.2710976651253826,0.2710976651253826)

# Row 6: sCs -> FAPs
Set-Row 6 @("sCs","Ccl25","Ccr10","FAPs",3,1,3.240161666666667,9.720485,0.2307926686810791,0.2307926686810791,2,0.6666666666666666,1.036595333333333,3.109786,0.393072250513715,0.393072250513715,3.358736462912222,30.22862816621,0.09071819368053795,0.09071819368053793)

# Row 7: sCs -> sCs
Set-Row 7 @("sCs","Ccl25","Ccr10","sCs",3,1,3.240161666666667,9.720485,0.2307926686810791,0.2307926686810791,3,1,1.600567,4.801701,0.6069277494862849,0.6069277494862849,5.186095838331667,46.674862544985,0.1400744750005411,0.1400744750005411)

Write-Output "applied edits"
